$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 302, shifting existing rows 302..372 down to 303..373
$ws.Rows("302").Insert()

# Populate the newly inserted row 302 with the new record's data.
# Columns A,B,C,E,F,G,H,N,O,Q,R stay constant across this block of rows.
$ws.Range("A302").Value = 9
$ws.Range("B302").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C302").Value = "Metropolitana"
$ws.Range("D302").Value = 44508
$ws.Range("E302").Value = 13
$ws.Range("F302").Value = 100112009
$ws.Range("G302").Value = "Acelga"
$ws.Range("H302").Value = "Sin especificar"
$ws.Range("I302").Value = "Primera"
$ws.Range("J302").Value = 52
$ws.Range("K302").Value = 10000
$ws.Range("L302").Value = 10000
$ws.Range("M302").Value = 10000
$ws.Range("N302").Value = "`$/docena de atados"
$ws.Range("O302").Value = "Región Metropolitana"
$ws.Range("P302").Value = 3333
$ws.Range("Q302").Value = 3
$ws.Range("R302").Value = "Hortaliza"
